$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 95
$ws1.Range("F4").Value = 73
$ws1.Range("F5").Value = 2514
$ws1.Range("F6").Value = 231
$ws1.Range("F7").Value = 379

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 95
$ws4.Range("F4").Value = 73
$ws4.Range("F5").Value = 2514
$ws4.Range("F6").Value = 231
$ws4.Range("F9").Value = 379

$wb.Save()
